$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update weekly triaged issues
$ws.Range("B2").Value = 63
$ws.Range("C2").Value = 3

$ws.Range("B12").Value = 30
$ws.Range("C12").Value = 21

$ws.Range("B13").Value = 9
$ws.Range("C13").Value = 10
